$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 550, shifting the existing weekly blocks
# (previously rows 550:581) down to 553:584, to make room for the newest
# week of "Acelga" price data at Mercado Mayorista Lo Valledor de Santiago.
$ws.Rows("550:552").Insert()

# Row 550: new weekly record - Extra quality
$ws.Range("A550").Value = 6
$ws.Range("B550").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C550").Value = "Metropolitana"
$ws.Range("D550").Value = 44516
$ws.Range("E550").Value = 13
$ws.Range("F550").Value = 100112009
$ws.Range("G550").Value = "Acelga"
$ws.Range("H550").Value = "Sin especificar"
$ws.Range("I550").Value = "Extra"
$ws.Range("J550").Value = 110
$ws.Range("K550").Value = 13000
$ws.Range("L550").Value = 13000
$ws.Range("M550").Value = 13000
$ws.Range("N550").Value = "`$/docena de atados"
$ws.Range("O550").Value = "Región Metropolitana"
$ws.Range("P550").Value = 4333
$ws.Range("Q550").Value = 3
$ws.Range("R550").Value = "Hortaliza"

# Row 551: new weekly record - Primera quality
$ws.Range("A551").Value = 6
$ws.Range("B551").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C551").Value = "Metropolitana"
$ws.Range("D551").Value = 44516
$ws.Range("E551").Value = 13
$ws.Range("F551").Value = 100112009
$ws.Range("G551").Value = "Acelga"
$ws.Range("H551").Value = "Sin especificar"
$ws.Range("I551").Value = "Primera"
$ws.Range("J551").Value = 150
$ws.Range("K551").Value = 11000
$ws.Range("L551").Value = 11000
$ws.Range("M551").Value = 11000
$ws.Range("N551").Value = "`$/docena de atados"
$ws.Range("O551").Value = "Región Metropolitana"
$ws.Range("P551").Value = 3667
$ws.Range("Q551").Value = 3
$ws.Range("R551").Value = "Hortaliza"

# Row 552: new weekly record - Segunda quality
$ws.Range("A552").Value = 6
$ws.Range("B552").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C552").Value = "Metropolitana"
$ws.Range("D552").Value = 44516
$ws.Range("E552").Value = 13
$ws.Range("F552").Value = 100112009
$ws.Range("G552").Value = "Acelga"
$ws.Range("H552").Value = "Sin especificar"
$ws.Range("I552").Value = "Segunda"
$ws.Range("J552").Value = 80
$ws.Range("K552").Value = 10000
$ws.Range("L552").Value = 10000
$ws.Range("M552").Value = 10000
$ws.Range("N552").Value = "`$/docena de atados"
$ws.Range("O552").Value = "Región Metropolitana"
$ws.Range("P552").Value = 3333
$ws.Range("Q552").Value = 3
$ws.Range("R552").Value = "Hortaliza"
